# Abrir orden de análisis
# Applies the changes described by the diff:
#  - Updates the selected/visible range of the sheet view
#  - Row 15: changes C15 text, fills in F15/G15/H15 (new row height 38.25)
#  - Row 16: fills in F16/G16/H16 (new row height 38.25)
#  - Row 17: fills in F17/G17/H17
#  - Rows 20-32: clears the "Nº CASO de PRUEBA" (column D) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 ---
$ws.Rows.Item(15).RowHeight = 38.25
$ws.Range("C15").Value = "Muestra para análisis"
$ws.Range("F15").Value = "Crear muestra de análisis"
$ws.Range("G15").Value = "Seleccionar atención e ingresar campo de comentario"
$ws.Range("H15").Value = "Creación de resultado de atención"

# --- Row 16 ---
$ws.Rows.Item(16).RowHeight = 38.25
$ws.Range("F16").Value = "Abrir una orden de análisis"
$ws.Range("G16").Value = "Seleccionar un resultado de atención y abrir la orden"
$ws.Range("H16").Value = "Abrir orden de análisis"

# --- Row 17 ---
$ws.Range("F17").Value = "Abrir orden sin seleccionar"
$ws.Range("G17").Value = "No seleccionar atención e intentar agendar"
$ws.Range("H17").Value = "Botón de abrir orden desactivado"

# --- Clear "Nº CASO de PRUEBA" numbering in column D for rows 20-32 ---
for ($r = 20; $r -le 32; $r++) {
    $ws.Range("D$r").Value = ""
}

# --- Update the view: scroll so row 14 is at the top, select C25 ---
$ws.Range("A14").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("C25").Select() | Out-Null
